# Improved Notifications and minnor bug fixes
#
# 1) ExamenesCalendario: append 7 new exam rows (rows 7-13)
# 2) ResumenEstadisticas: refresh the "generated at" timestamp and the
#    exam-related counters that change because of the new rows above.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) ExamenesCalendario - append new rows
# ---------------------------------------------------------------------
$wsCal = $wb.Worksheets.Item("ExamenesCalendario")

$newRows = @(
    @{ Date = 45820; Turn = "1"; Type = "suficiencia"; Subject = "b"   },
    @{ Date = 45820; Turn = "1"; Type = "premio";      Subject = "cv"  },
    @{ Date = 45820; Turn = "1"; Type = "suficiencia"; Subject = "asd" },
    @{ Date = 45820; Turn = "1"; Type = "suficiencia"; Subject = "fs"  },
    @{ Date = 45820; Turn = "2"; Type = "suficiencia"; Subject = "nom" },
    @{ Date = 45821; Turn = "1"; Type = "premio";      Subject = "sd2" },
    @{ Date = 45821; Turn = "1"; Type = "suficiencia"; Subject = "grt" }
)

$startRow = 7
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    # Date column keeps the same date number format used by the existing
    # rows in column A.
    $wsCal.Cells.Item($r, 1).Value = $row.Date
    $wsCal.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD"

    # Turn column holds text like "1"/"2" - force text formatting so the
    # numeric-looking value isn't coerced into a real number, then clear
    # the formatting again so no extra style gets attached to the cell
    # (matching the rest of the column which has no explicit style).
    $wsCal.Cells.Item($r, 2).NumberFormat = "@"
    $wsCal.Cells.Item($r, 2).Value = $row.Turn
    $wsCal.Cells.Item($r, 2).ClearFormats()

    $wsCal.Cells.Item($r, 3).Value = $row.Type
    $wsCal.Cells.Item($r, 4).Value = $row.Subject
}

# ---------------------------------------------------------------------
# 2) ResumenEstadisticas - refresh report timestamp & counters
# ---------------------------------------------------------------------
$wsStats = $wb.Worksheets.Item("ResumenEstadisticas")

$wsStats.Cells.Item(2, 2).Value = "2025-06-09 14:16:57"  # Fecha de Generación del Reporte
$wsStats.Cells.Item(7, 2).Value = 12                      # Total de Exámenes en Calendario
$wsStats.Cells.Item(8, 2).Value = 7                       # Exámenes por Tipo - Suficiencia
$wsStats.Cells.Item(9, 2).Value = 5                       # Exámenes por Tipo - Premio
